# Apply cryptos.xlsx updates: refreshed prices/volumes and a few reordered rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $row, $col, $val) {
    # Force the cell to stay a text string even if the value looks numeric,
    # then restore the default style so no stray formatting is introduced.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '27.932.11'
$ws.Cells.Item(2, 5).Value = '  +1.52%  '

$ws.Cells.Item(3, 4).Value = '1.641.61'
$ws.Cells.Item(3, 5).Value = '  +1.25%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

Set-TextValue $ws 5 4 '213.58'
$ws.Cells.Item(5, 5).Value = '  +1.03%  '

$ws.Cells.Item(6, 5).Value = '  +0.17%  '

$ws.Cells.Item(7, 5).Value = '  -0.07%  '

Set-TextValue $ws 8 4 '23.76'
$ws.Cells.Item(8, 5).Value = '  +3.19%  '

$ws.Cells.Item(9, 5).Value = '  +0.32%  '

$ws.Cells.Item(10, 5).Value = '  +0.89%  '

$ws.Cells.Item(11, 5).Value = '  -0.56%  '

$ws.Cells.Item(12, 4).Value = '1.874.02'
$ws.Cells.Item(12, 5).Value = '  +1.26%  '

$ws.Cells.Item(13, 4).Value = '1.640.54'
$ws.Cells.Item(13, 5).Value = '  +0.95%  '

Set-TextValue $ws 15 4 '0.574'
$ws.Cells.Item(15, 5).Value = '  +4.48%  '

Set-TextValue $ws 16 4 '66.15'
$ws.Cells.Item(16, 5).Value = '  +1.46%  '

$ws.Cells.Item(17, 4).Value = '27.919.00'
$ws.Cells.Item(17, 5).Value = '  +1.57%  '

Set-TextValue $ws 18 4 '231.20'
$ws.Cells.Item(18, 5).Value = '  +0.74%  '

$ws.Cells.Item(19, 5).Value = '  +1.30%  '

Set-TextValue $ws 20 4 '7.62'

Set-TextValue $ws 21 4 '11.38'
$ws.Cells.Item(21, 5).Value = '  +9.62%  '

$ws.Cells.Item(22, 5).Value = '  -0.06%  '

$ws.Cells.Item(23, 5).Value = '  +1.54%  '

$ws.Cells.Item(24, 5).Value = '  -2.28%  '

Set-TextValue $ws 25 4 '152.14'
$ws.Cells.Item(25, 5).Value = '  +1.89%  '

Set-TextValue $ws 26 4 '6.94'

$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 27 4 '15.74'
$ws.Cells.Item(27, 5).Value = '  +1.34%  '

$ws.Cells.Item(28, 2).Value = 'Stellar'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 28 4 '0.112'
$ws.Cells.Item(28, 5).Value = '  +0.77%  '

$ws.Cells.Item(30, 5).Value = '  +1.08%  '

Set-TextValue $ws 31 4 '0.0486'
$ws.Cells.Item(31, 5).Value = '  +0.65%  '

$ws.Cells.Item(32, 5).Value = '  +1.97%  '

$ws.Cells.Item(33, 4).Value = '1.422.93'
$ws.Cells.Item(33, 5).Value = '  -2.78%  '

Set-TextValue $ws 34 4 '3.12'
$ws.Cells.Item(34, 5).Value = '  +2.32%  '

$ws.Cells.Item(35, 5).Value = '  +1.84%  '

Set-TextValue $ws 36 4 '2.35'
$ws.Cells.Item(36, 5).Value = '  +0.38%  '

Set-TextValue $ws 37 4 '0.892'
$ws.Cells.Item(37, 5).Value = '  +2.46%  '

$ws.Cells.Item(38, 5).Value = '  +0.70%  '

Set-TextValue $ws 39 4 '0.923'
$ws.Cells.Item(39, 5).Value = '  -1.96%  '

$ws.Cells.Item(40, 5).Value = '  +1.26%  '

$ws.Cells.Item(41, 5).Value = '  +2.32%  '

Set-TextValue $ws 43 4 '67.13'
$ws.Cells.Item(43, 5).Value = '  -0.24%  '

$ws.Cells.Item(44, 5).Value = '  +0.58%  '

$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 45 4 '1.82'
$ws.Cells.Item(45, 5).Value = '  +3.93%  '

$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 46 4 '5.45'
$ws.Cells.Item(46, 5).Value = '  +2.95%  '

$ws.Cells.Item(47, 5).Value = '  +0.19%  '

$ws.Cells.Item(48, 4).Value = '1.782.72'
$ws.Cells.Item(48, 5).Value = '  +1.24%  '

Set-TextValue $ws 49 4 '88.79'
$ws.Cells.Item(49, 5).Value = '  +1.91%  '

$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 50 4 '0.101'
$ws.Cells.Item(50, 5).Value = '  +1.38%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 51 4 '0.0506'
$ws.Cells.Item(51, 5).Value = '  +0.63%  '
